# -------------------------------------------------------------------------
# Adds "Sheet2" (calcLists checks) after Sheet1, populates it with the
# checklist / averaging-example content, restyles the header rows and
# restores the selection/zoom state described in the commit.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet1: selection changes to B6:C38 (tabSelected moves to Sheet2) ----
$ws1.Range("B6:C38").Select() | Out-Null

# ---- Create Sheet2 right after Sheet1 ----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# ---- Column widths (approximate character widths used in the source) ----
$ws2.Columns.Item(1).ColumnWidth = 9.166666666666666
$ws2.Columns.Item(2).ColumnWidth = 40.166666666666664
$ws2.Columns.Item(3).ColumnWidth = 33.5
$ws2.Columns.Item(4).ColumnWidth = 7.833333333333333
$ws2.Columns.Item(6).ColumnWidth = 91.83333333333334
$ws2.Columns.Item(7).ColumnWidth = 11.833333333333332
$ws2.Columns.Item(11).ColumnWidth = 15.0

# ---- Row 1: checks header ----
$ws2.Range("A1").Value = "Region & own effect"
$ws2.Range("C1").Value = "Checks"
$ws2.Range("D1").Value = "Row"
$ws2.Range("F1").Value = "Notes"
$ws2.Range("A1,C1,D1,F1,G1").Font.Bold = $true
$ws2.Range("D1").HorizontalAlignment = -4152

# ---- Row 2-3: US national ----
$ws2.Range("B2").Value = "US national [as now]"
$ws2.Range("B2").Font.Bold = $true
$ws2.Range("C2").Value = "US National"
$ws2.Range("D2").Value = 28
$ws2.Range("C3").Value = "NO Cross Effects (cross-effect = 0)"
$ws2.Range("D3").Value = 81

# ---- Row 5-6: Corn belt ----
$ws2.Range("B5").Value = "Corn belt (all or some part) " + [char]0x2013 + " must add one row"
$ws2.Range("B5").Font.Bold = $true
$ws2.Range("C5").Value = "Corn Belt"
$ws2.Range("D5").Value = 30
$ws2.Range("C6").Value = "NO Cross Effects (cross-effect = 0)"
$ws2.Range("D6").Value = 81

# ---- Row 8-9: Other ----
$ws2.Range("B8").Value = "Other [as now]"
$ws2.Range("B8").Font.Bold = $true
$ws2.Range("C8").Value = "Everything except US and Corn Belt"
$ws2.Range("D8").Value = "t-(28 + 30)"
$ws2.Range("C9").Value = "NO Cross Effects (cross-effect = 0)"
$ws2.Range("D9").Value = 81

# ---- Row 11: Methods & own effect ----
$ws2.Range("A11").Value = "Methods & own effect"
$ws2.Range("A11").Font.Bold = $true

# ---- Row 12-15: Estimation, panel or survey data ----
$ws2.Range("B12").Value = "Estimation, panel or survey data"
$ws2.Range("B12").Font.Bold = $true
$ws2.Range("C12").Value = "Estimated: balanced panel data"
$ws2.Range("D12").Value = 20
$ws2.Range("C13").Value = "Estimated: unbalanced panel data"
$ws2.Range("D13").Value = 21
$ws2.Range("C14").Value = "Estimated: survey"
$ws2.Range("D14").Value = 19
$ws2.Range("F14").Value = "Do we also need to add non estimate survey row? (24)"
$ws2.Range("C15").Value = "NO Cross Effects (cross-effect = 0)"
$ws2.Range("D15").Value = 81

# ---- Row 17-18: Estimation, market data ----
$ws2.Range("B17").Value = "Estimation, market data"
$ws2.Range("B17").Font.Bold = $true
$ws2.Range("C17").Value = "Estimated: market data"
$ws2.Range("D17").Value = 18
$ws2.Range("C18").Value = "NO Cross Effects (cross-effect = 0)"
$ws2.Range("D18").Value = 81

# ---- Row 20-22: Simulation or theory ----
$ws2.Range("B20").Value = "Simulation or theory"
$ws2.Range("B20").Font.Bold = $true
$ws2.Range("C20").Value = "Simulation"
$ws2.Range("D20").Value = 22
$ws2.Range("C21").Value = "Theory"
$ws2.Range("D21").Value = 23
$ws2.Range("C22").Value = "NO Cross Effects (cross-effect = 0)"
$ws2.Range("D22").Value = 81

# ---- Row 24: All & own effect ----
$ws2.Range("A24").Value = "All & own effect"
$ws2.Range("A24").Font.Bold = $true
$ws2.Range("C24").Value = "NO Cross Effects (cross-effect = 0)"
$ws2.Range("D24").Value = 81

# ---- Row 27-28: Nature of effect / Cross-effect ----
$ws2.Range("A27").Value = "Nature of effect"
$ws2.Range("A27").Font.Bold = $true
$ws2.Range("B28").Value = "Cross-effect"
$ws2.Range("B28").Font.Bold = $true
$ws2.Range("C28").Value = "Cross Effects (cross-effect = 1)"
$ws2.Range("D28").Value = 81
$ws2.Range("F28").Value = """Opposite"" of all and own effect"

# ---- Row 31-32: All crops ----
$ws2.Range("B31").Value = "All crops"
$ws2.Range("B31").Font.Bold = $true
$ws2.Range("C31").Value = "area of all crops; production of all crops"
$ws2.Range("D31").Value = "71, 74"
$ws2.Range("F31").Value = "Dependent on Area, yield, or production"
$ws2.Range("F32").Value = "Yield = N/a"

# ---- Row 34: One crop ----
$ws2.Range("B34").Value = "One crop"
$ws2.Range("B34").Font.Bold = $true
$ws2.Range("C34").Value = "area of one crop, production of one crop"
$ws2.Range("D34").Value = "70, 73"
$ws2.Range("D34").HorizontalAlignment = -4152
$ws2.Range("F34").Value = "Dependent on Area, yield, or production"

# ---- Highlight the "checks" columns C & D (yellow fill) ----
$yellow = 65535
$ws2.Range("C2:D3").Interior.Color = $yellow
$ws2.Range("C5:D6").Interior.Color = $yellow
$ws2.Range("C8:D9").Interior.Color = $yellow
$ws2.Range("C12:D15").Interior.Color = $yellow
$ws2.Range("C17:D18").Interior.Color = $yellow
$ws2.Range("C20:D22").Interior.Color = $yellow
$ws2.Range("C24:D24").Interior.Color = $yellow
$ws2.Range("C28:D28").Interior.Color = $yellow
$ws2.Range("C31:D31").Interior.Color = $yellow

# Right-align the numeric/"row" check column (D) for every highlighted block
$ws2.Range("D2:D3").HorizontalAlignment = -4152
$ws2.Range("D5:D6").HorizontalAlignment = -4152
$ws2.Range("D8:D9").HorizontalAlignment = -4152
$ws2.Range("D12:D15").HorizontalAlignment = -4152
$ws2.Range("D17:D18").HorizontalAlignment = -4152
$ws2.Range("D20:D22").HorizontalAlignment = -4152
$ws2.Range("D24").HorizontalAlignment = -4152
$ws2.Range("D28").HorizontalAlignment = -4152
$ws2.Range("D31").HorizontalAlignment = -4152

# -------------------------------------------------------------------------
# Weighted-average example (rows 40-50)
# -------------------------------------------------------------------------
$ws2.Range("G40").Value = "Study"
$ws2.Range("H40").Value = "Value"
$ws2.Range("I40").Value = "Average"
$ws2.Range("J40").Value = "Weight"
$ws2.Range("K40").Value = "Weighted average"

$ws2.Range("F41").Value = "study-weighted average would give each study equal weight no matter how many observations each one has"
$ws2.Range("F41").Font.Bold = $true
$ws2.Range("G41").Value = 1
$ws2.Range("H41").Value = 4

$ws2.Range("F42").Value = "Example:"
$ws2.Range("G42").Value = 1
$ws2.Range("H42").Value = 7

$ws2.Range("F43").Value = "3 studies are included in the average. (1,2,3,4) (5,6) (7,8,9) = simple avg 5"
$ws2.Range("G43").Value = 1
$ws2.Range("H43").Value = 2

$ws2.Range("F44").Value = "weighted average: (1,2,3,4)/4 + (5,6)/2 + (7,8,9)/3 = weighted average 5.33"
$ws2.Range("G44").Value = 1
$ws2.Range("H44").Value = 7
$ws2.Range("I44").Formula = "=AVERAGE(H41:H44)"
$ws2.Range("J44").Value = "(1/3)"
$ws2.Range("J44").NumberFormat = "d-mmm"
$ws2.Range("K44").Formula = "=I44*(1/3)"

$ws2.Range("G45").Value = 2
$ws2.Range("H45").Value = 1

$ws2.Range("G46").Value = 2
$ws2.Range("H46").Value = 8
$ws2.Range("I46").Formula = "=AVERAGE(H45:H46)"
$ws2.Range("J46").Value = "(1/3)"
$ws2.Range("J46").NumberFormat = "d-mmm"
$ws2.Range("K46").Formula = "=I46*(1/3)"

$ws2.Range("G47").Value = 3
$ws2.Range("H47").Value = 9

$ws2.Range("G48").Value = 3
$ws2.Range("H48").Value = 7

$ws2.Range("G49").Value = 3
$ws2.Range("H49").Value = 3
$ws2.Range("I49").Formula = "=AVERAGE(H47:H49)"
$ws2.Range("J49").Value = "(1/3)"
$ws2.Range("J49").NumberFormat = "d-mmm"
$ws2.Range("K49").Formula = "=I49*(1/3)"

$ws2.Range("I50").Formula = "=AVERAGE(I41:I49)"
$ws2.Range("K50").Formula = "=SUM(K44:K49)"
$ws2.Range("K50").Font.Bold = $true

# Right-align the "study value" column (H)
$ws2.Range("H41:H49").HorizontalAlignment = -4152

# ---- Sheet2 view: zoom 70%, scrolled to row 19, selection C31:D31, active ----
$ws2.Activate()
$excel.ActiveWindow.Zoom = 70
$ws2.Range("C31:D31").Select() | Out-Null

Write-Host "Sheet2 created and populated."
